# Update the "想去人数" (want-to-go count) figures that were refreshed
# by the data scrape, on both the "展览" sheet and the "全部类型" sheet.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions) - column F holds the want-to-go counts.
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    "F3"  = 24
    "F4"  = 16233
    "F5"  = 422
    "F8"  = 15568
    "F10" = 9204
    "F11" = 454
    "F14" = 118
    "F17" = 218
    "F19" = 82
    "F20" = 597
    "F23" = 73
    "F24" = 1141
    "F27" = 32
    "F28" = 517
    "F30" = 45
    "F32" = 79
    "F36" = 359
    "F37" = 474
    "F39" = 5659
}
foreach ($addr in $expoUpdates.Keys) {
    $wsExpo.Range($addr).Value = $expoUpdates[$addr]
}

# Sheet "全部类型" (all types) - same underlying rows, but row numbers are
# shifted a bit further down due to extra rows in this consolidated sheet.
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    "F3"  = 24
    "F4"  = 16233
    "F5"  = 422
    "F8"  = 15568
    "F10" = 9204
    "F11" = 454
    "F14" = 118
    "F17" = 218
    "F19" = 82
    "F20" = 597
    "F23" = 73
    "F24" = 1141
    "F27" = 32
    "F28" = 517
    "F30" = 45
    "F34" = 79
    "F38" = 359
    "F39" = 474
    "F41" = 5659
}
foreach ($addr in $allUpdates.Keys) {
    $wsAll.Range($addr).Value = $allUpdates[$addr]
}
